$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Egress": the last data row (row 34) had two identical "Reserved"
# cells (E34 and F34). They get merged into a single E34:F34 cell, so F34's
# separate value is dropped and a new merge entry is added.
# ---------------------------------------------------------------------------
$egress = $wb.Worksheets.Item("Egress")
$egress.Range("F34").Clear()
$egress.Range("E34:F34").Merge()

# ---------------------------------------------------------------------------
# Sheet "Ingress": a new "Statuses" row is inserted before the final
# "Reserved" row, and the stale per-row "Reserved" placeholders in columns
# E/F collapse into one merged cell running down the whole data block.
# ---------------------------------------------------------------------------
$ingress = $wb.Worksheets.Item("Ingress")

# Insert a new row at 7 - this pushes the old row 7 ("Reserved" / "4-31")
# down to row 8, and Excel auto-extends merges that spanned row 7 (the B
# and C column merges, plus the stray E/F merges further down the sheet).
$ingress.Rows.Item(7).Insert()

# The B-column merge auto-grew to B3:B8; split it back apart since the new
# "Statuses" row (7) and the "Reserved" row (8) each carry their own B
# value now, separate from the B3:B6 "N/A" block.
$ingress.Range("B3:B8").UnMerge()

# Fill in the new "Statuses" row.
$ingress.Range("A7").Value = "Statuses"
$ingress.Range("B7").Value = 0.1
$ingress.Range("D7").Value = 4.0

# Match B7's style to the right-aligned look used by the other "value"
# column cells (copy formatting only, so the existing style is reused).
$ingress.Range("D8").Copy() | Out-Null
$ingress.Range("B7").PasteSpecial(-4122) | Out-Null

# Former row 7 (now row 8): add its own "N/A" in column B, and update the
# stale "4-31" bit-range label to "5-31".
$ingress.Range("B8").Value = "N/A"
$ingress.Range("D8").Copy() | Out-Null
$ingress.Range("B8").PasteSpecial(-4122) | Out-Null
$ingress.Range("D8").Value = "5-31"

# Re-merge the B column for just the unchanged "N/A" block (rows 3-6).
$ingress.Range("B3:B6").Merge()

# The per-row "Reserved" labels in columns E/F (rows 3-8) collapse into one
# merged cell, keeping only the row-3 value.
$ingress.Range("F3").Clear()
$ingress.Range("E4:F6").Clear()
$ingress.Range("E7:F8").Clear()
$ingress.Range("E3:F8").Merge()
